$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as text so values like "311.70" or "0.630"
# are not silently re-interpreted/rounded as numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.990.33"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.356.88"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.70"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.19"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.79"
$ws.Range("E10").Value = "  -2.91%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.43"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.971"
$ws.Range("E14").Value = "  -4.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.715.82"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.18"
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.362.17"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.980.28"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.00"
$ws.Range("E19").Value = "  +7.49%  "
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.16"
$ws.Range("E21").Value = "  -5.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.01"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "258.64"
$ws.Range("E24").Value = "  -3.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.30"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.00"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.16"
$ws.Range("E28").Value = "  -6.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  -3.36%  "
$ws.Range("E30").Value = "  +5.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.24"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.00"
$ws.Range("E32").Value = "  -6.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "167.27"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.00"
$ws.Range("E34").Value = "  +4.16%  "
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0351"
$ws.Range("E39").Value = "  -4.39%  "
$ws.Range("E40").Value = "  -1.81%  "
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.04"
$ws.Range("E42").Value = "  -5.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.32"
$ws.Range("E44").Value = "  -5.64%  "
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.842.10"
$ws.Range("E46").Value = "  +10.91%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.75"
$ws.Range("E47").Value = "  -8.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "83.22"
$ws.Range("E48").Value = "  +6.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.67"
$ws.Range("E49").Value = "  +6.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.42"
$ws.Range("E50").Value = "  -4.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.14"
$ws.Range("E51").Value = "  +1.11%  "
